# 货物服务消费占比.xlsx update:
#   - swap the (A,B,C) data for the "B"/"C" quarter rows within each year
#     group (2019, 2020, 2021) back to their correct order
#   - drop the now-redundant D/E columns (they duplicated B/C with some
#     stray artifacts)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Row($r1, $r2) {
    $a1 = $ws.Range("A$r1").Value2
    $b1 = $ws.Range("B$r1").Value2
    $c1 = $ws.Range("C$r1").Value2
    $a2 = $ws.Range("A$r2").Value2
    $b2 = $ws.Range("B$r2").Value2
    $c2 = $ws.Range("C$r2").Value2

    $ws.Range("A$r1").Value = $a2
    $ws.Range("B$r1").Value = $b2
    $ws.Range("C$r1").Value = $c2

    $ws.Range("A$r2").Value = $a1
    $ws.Range("B$r2").Value = $b1
    $ws.Range("C$r2").Value = $c1
}

Swap-Row 5 6
Swap-Row 9 10
Swap-Row 13 14

# Remove columns D and E entirely (header + all data rows)
$ws.Range("D1:E15").Delete()
